# Apply updated crypto price / 1h-volume data (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''25.755.66'
$ws.Range('E2').Value = '  -0.93%  '
$ws.Range('D3').Value = '''1.624.20'
$ws.Range('E3').Value = '  -0.88%  '
$ws.Range('D4').Value = '''1.002'
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').Value = '''214.80'
$ws.Range('E5').Value = '  -0.09%  '
$ws.Range('D6').Value = '''0.5067'
$ws.Range('E6').Value = '  -1.03%  '
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('D8').Value = '''0.2558'
$ws.Range('E8').Value = '  -0.91%  '
$ws.Range('D9').Value = '''0.06377'
$ws.Range('E9').Value = '  +0.18%  '
$ws.Range('D10').Value = '''19.27'
$ws.Range('E10').Value = '  -2.71%  '
$ws.Range('D11').Value = '''0.07770'
$ws.Range('E12').Value = '  -1.01%  '
$ws.Range('D13').Value = '''1.623.72'
$ws.Range('E13').Value = '  -1.78%  '
$ws.Range('D14').Value = '''1.847.38'
$ws.Range('E14').Value = '  -1.09%  '
$ws.Range('D15').Value = '''0.5538'
$ws.Range('E15').Value = '  +1.34%  '
$ws.Range('D16').Value = '''63.55'
$ws.Range('E16').Value = '  -1.64%  '
$ws.Range('D17').Value = '''0.0₅7511'
$ws.Range('E17').Value = '  -2.97%  '
$ws.Range('D18').Value = '''25.774.57'
$ws.Range('E18').Value = '  -0.90%  '
$ws.Range('D19').Value = '''1.001'
$ws.Range('E19').Value = '  -0.16%  '
$ws.Range('D20').Value = '''193.56'
$ws.Range('E20').Value = '  -1.99%  '
$ws.Range('D21').Value = '''4.389'
$ws.Range('E21').Value = '  -1.05%  '
$ws.Range('D22').Value = '''9.750'
$ws.Range('E22').Value = '  -2.24%  '
$ws.Range('D23').Value = '''5.963'
$ws.Range('E23').Value = '  -2.02%  '
$ws.Range('E24').Value = '  -0.29%  '
$ws.Range('D25').Value = '''1.863'
$ws.Range('E25').Value = '  -1.56%  '
$ws.Range('D26').Value = '''140.69'
$ws.Range('E26').Value = '  -0.91%  '
$ws.Range('E27').Value = '  +1.02%  '
$ws.Range('E28').Value = '  -1.93%  '
$ws.Range('D29').Value = '''15.42'
$ws.Range('E29').Value = '  -1.53%  '
$ws.Range('E30').Value = '  -0.39%  '
$ws.Range('D31').Value = '''0.04855'
$ws.Range('E31').Value = '  -0.75%  '
$ws.Range('D32').Value = '''3.306'
$ws.Range('E32').Value = '  +0.76%  '
$ws.Range('D33').Value = '''3.172'
$ws.Range('E33').Value = '  -1.24%  '
$ws.Range('E34').Value = '  +0.18%  '
$ws.Range('E35').Value = '  -0.60%  '
$ws.Range('E36').Value = '  -2.54%  '
$ws.Range('D37').Value = '''1.124.36'
$ws.Range('E37').Value = '  +1.14%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').Value = '''2.535'
$ws.Range('E38').Value = '  -1.99%  '
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').Value = '''0.5491'
$ws.Range('E39').Value = '  -0.86%  '
$ws.Range('D40').Value = '''0.01555'
$ws.Range('E41').Value = '  -0.19%  '
$ws.Range('D42').Value = '''5.560'
$ws.Range('E42').Value = '  +0.51%  '
$ws.Range('D43').Value = '''0.7931'
$ws.Range('E43').Value = '  -1.98%  '
$ws.Range('D44').Value = '''97.08'
$ws.Range('E44').Value = '  -2.37%  '
$ws.Range('D45').Value = '''1.770.60'
$ws.Range('E45').Value = '  -0.34%  '
$ws.Range('D46').Value = '''0.0₈115'
$ws.Range('E46').Value = '  -6.04%  '
$ws.Range('D47').Value = '''0.4418'
$ws.Range('E47').Value = '  -2.59%  '
$ws.Range('D48').Value = '''54.61'
$ws.Range('E48').Value = '  -0.82%  '
$ws.Range('D49').Value = '''0.05129'
$ws.Range('E49').Value = '  -3.18%  '
$ws.Range('D50').Value = '''7.559'
$ws.Range('E50').Value = '  +2.83%  '
$ws.Range('D51').Value = '''0.9989'
$ws.Range('E51').Value = '  -0.80%  '
